$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: 889×8= -> 735×3=
$cell = $t.Cell(1, 1)
$cell.Range.Text = "735×3="

# Row 1, Col 2: 747×4= -> 794×2=
$cell = $t.Cell(1, 2)
$cell.Range.Text = "794×2="

# Row 1, Col 3: 776×3= -> 919×4=
$cell = $t.Cell(1, 3)
$cell.Range.Text = "919×4="

# Row 1, Col 4: 720×8= -> 545×3=
$cell = $t.Cell(1, 4)
$cell.Range.Text = "545×3="

# Row 1, Col 5: 302×3= -> 480×6=
$cell = $t.Cell(1, 5)
$cell.Range.Text = "480×6="

# Row 5, Col 1: 454×6= -> 480×6=
$cell = $t.Cell(5, 1)
$cell.Range.Text = "480×6="

# Row 5, Col 2: 965×5= -> 749×3=
$cell = $t.Cell(5, 2)
$cell.Range.Text = "749×3="

# Row 5, Col 3: 184×6= -> 846×6=
$cell = $t.Cell(5, 3)
$cell.Range.Text = "846×6="

# Row 5, Col 4: 819×6= -> 574×4=
$cell = $t.Cell(5, 4)
$cell.Range.Text = "574×4="

# Row 5, Col 5: 895×9= -> 482×2=
$cell = $t.Cell(5, 5)
$cell.Range.Text = "482×2="

# Row 10, Col 1: 754×3= -> 719×4=
$cell = $t.Cell(10, 1)
$cell.Range.Text = "719×4="

# Row 10, Col 2: 589×8= -> 806×3=
$cell = $t.Cell(10, 2)
$cell.Range.Text = "806×3="

# Row 10, Col 3: 200×6= -> 776×2=
$cell = $t.Cell(10, 3)
$cell.Range.Text = "776×2="

# Row 10, Col 4: 321×2= -> 747×6=
$cell = $t.Cell(10, 4)
$cell.Range.Text = "747×6="

# Row 10, Col 5: 361×5= -> 993×9=
$cell = $t.Cell(10, 5)
$cell.Range.Text = "993×9="

# Row 15, Col 1: 280×6= -> 354×8=
$cell = $t.Cell(15, 1)
$cell.Range.Text = "354×8="

# Row 15, Col 2: 490×3= -> 526×9=
$cell = $t.Cell(15, 2)
$cell.Range.Text = "526×9="

# Row 15, Col 3: 677×4= -> 931×3=
$cell = $t.Cell(15, 3)
$cell.Range.Text = "931×3="

# Row 15, Col 4: 749×3= -> 337×3=
$cell = $t.Cell(15, 4)
$cell.Range.Text = "337×3="

# Row 15, Col 5: 424×2= -> 588×9=
$cell = $t.Cell(15, 5)
$cell.Range.Text = "588×9="

# Row 20, Col 1: 306×6= -> 642×4=
$cell = $t.Cell(20, 1)
$cell.Range.Text = "642×4="

# Row 20, Col 2: 143×2= -> 101×6=
$cell = $t.Cell(20, 2)
$cell.Range.Text = "101×6="

# Row 20, Col 3: 594×4= -> 837×3=
$cell = $t.Cell(20, 3)
$cell.Range.Text = "837×3="

# Row 20, Col 4: 175×9= -> 620×3=
$cell = $t.Cell(20, 4)
$cell.Range.Text = "620×3="

# Row 20, Col 5: 305×4= -> 515×8=
$cell = $t.Cell(20, 5)
$cell.Range.Text = "515×8="
